# This document's 8 "content" paragraphs get their text values rotated in a
# single cycle (each paragraph's new text is the OLD text of the paragraph
# that, in document order, comes right after it; the last one wraps to the
# second one):
#   Objetivos paragraph            <- old "Programa resumido" text
#   Docente paragraph               <- old "Objetivos" text
#   Programa resumido paragraph     <- old "Programa" text
#   Programa paragraph              <- old "Método" text
#   Método value                    <- old "Critério" value
#   Critério value                  <- old "Norma de recuperação" value
#   Norma de recuperação value      <- old "Bibliografia" text
#   Bibliografia paragraph          <- old "Docente" text
#
# We do the replacements with Find/Replace using distinct, never-colliding
# placeholder markers first (old text -> placeholder, at each location),
# then a second pass turns each location's placeholder into the real final
# text for that location. This avoids any ordering hazard where a later
# replacement would accidentally match text just written by an earlier one.

$d = $word.ActiveDocument

$txtObjetivos   = "Fornecer oportunidade de realização de treinamento profissional de Engenharia Ambiental em empresa ou instituição sob supervisão de docente do Departamento de Ciências Básicas e Ambientais da EEL. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
$txtDocente     = "4780627 - Ana Lucia Gabas Ferreira"
$txtProgResumido= "Processo seletivo. Plano de trabalho específico. Realização do estágio. Relatório final."
$txtPrograma    = "Participação do aluno em processo seletivo de empresas, instituições de pesquisa ou no setor acadêmico. O estágio realizado sob a supervisão de docente designado pelo Departamento de Ciências Básicas e Ambientais da Escola de Engenharia de Lorena. O conteúdo será estabelecido no Plano de Trabalho entre o supervisor responsável pelo Estágio e o docente supervisor. Apresentação de relatório final sobre as atividades desenvolvidas no estágio."
$txtMetodo      = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$txtCriterio    = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$txtNormaRecup  = "Devido às características da disciplina, não será oferecida recuperação."
$txtBibliografia= "Não há."

# Step 1: old text -> unique placeholder (so subsequent Find passes can't
# re-match text that a previous replacement just inserted).
$d.Content.Find.Execute($txtObjetivos, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_OBJETIVOS@@", 2)
$d.Content.Find.Execute($txtDocente, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_DOCENTE@@", 2)
$d.Content.Find.Execute($txtProgResumido, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_PROGRESUMIDO@@", 2)
$d.Content.Find.Execute($txtPrograma, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_PROGRAMA@@", 2)
$d.Content.Find.Execute($txtMetodo, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_METODO@@", 2)
$d.Content.Find.Execute($txtCriterio, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_CRITERIO@@", 2)
$d.Content.Find.Execute($txtNormaRecup, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_NORMARECUP@@", 2)
$d.Content.Find.Execute($txtBibliografia, $true, $false, $false, $false, $false, $true, 1, $false, "@@PH_BIBLIOGRAFIA@@", 2)

# Step 2: placeholder -> the new, rotated text content.
# Rule: the placeholder now sitting at a given paragraph's location gets
# replaced with whatever text the diff says that paragraph's NEW value is.
$d.Content.Find.Execute("@@PH_OBJETIVOS@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtProgResumido, 2)
$d.Content.Find.Execute("@@PH_DOCENTE@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtObjetivos, 2)
$d.Content.Find.Execute("@@PH_PROGRESUMIDO@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtPrograma, 2)
$d.Content.Find.Execute("@@PH_PROGRAMA@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtMetodo, 2)
$d.Content.Find.Execute("@@PH_METODO@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtCriterio, 2)
$d.Content.Find.Execute("@@PH_CRITERIO@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtNormaRecup, 2)
$d.Content.Find.Execute("@@PH_NORMARECUP@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtBibliografia, 2)
$d.Content.Find.Execute("@@PH_BIBLIOGRAFIA@@", $true, $false, $false, $false, $false, $true, 1, $false, $txtDocente, 2)
